$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 410.5
$ws.Range("I20").Value = 410.5
$ws.Range("K20").Value = 410.5
$ws.Range("M20").Value = -180.5

$ws.Range("H35").Value = 410.5
$ws.Range("I35").Value = 410.5
$ws.Range("K35").Value = 410.5
$ws.Range("M35").Value = -31.5

$ws.Range("H76").Value = 3777.6667

$ws.Range("H79").Value = 3777.6667

$ws.Range("H80").Value = 7858654.5
$ws.Range("I80").Value = 644.0
$ws.Range("J80").Value = 15225540.0
$ws.Range("K80").Value = 1932.0
$ws.Range("L80").Value = 45676620.0
$ws.Range("M80").Value = -934.0
$ws.Range("N80").Value = -45678616.0

$ws.Range("H83").Value = 7858654.5
$ws.Range("I83").Value = 644.0
$ws.Range("J83").Value = 15225540.0
$ws.Range("K83").Value = 5796.0
$ws.Range("L83").Value = 137029860.0
$ws.Range("M83").Value = -804.0
$ws.Range("N83").Value = -137039844.0

$ws.Range("H111").Value = 1964.7142
$ws.Range("I111").Value = 2125.5
$ws.Range("J111").Value = 1000.0
$ws.Range("K111").Value = 6376.5
$ws.Range("L111").Value = 3000.0
$ws.Range("M111").Value = -3309.5
$ws.Range("N111").Value = -9134.0

$ws.Range("H112").Value = 1079.3125
$ws.Range("J112").Value = 1079.3125
$ws.Range("L112").Value = 3237.9375
$ws.Range("N112").Value = -5453.9375

$ws.Range("H132").Value = 3993.5715
$ws.Range("I132").Value = 4509.1665
$ws.Range("K132").Value = 13527.4995
$ws.Range("M132").Value = -10997.4995

$ws.Range("H137").Value = 1894.6364
$ws.Range("I137").Value = 1476.0769
$ws.Range("K137").Value = 4428.2307
$ws.Range("M137").Value = -1878.2307

$ws.Range("H138").Value = 2096.494
$ws.Range("I138").Value = 1732.4445
$ws.Range("J138").Value = 2200.508
$ws.Range("K138").Value = 5197.333500000001
$ws.Range("L138").Value = 6601.523999999999
$ws.Range("M138").Value = -57.33350000000064
$ws.Range("N138").Value = -16881.524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1737.3715
$ws.Range("J2").Value = 2112.125
$ws.Range("L2").Value = 2112.125
$ws.Range("N2").Value = -2338.125

$ws.Range("H32").Value = 5137.293
$ws.Range("I32").Value = 5354.625
$ws.Range("K32").Value = 5354.625
$ws.Range("M32").Value = -5067.625

$ws.Range("H110").Value = 556.44446
$ws.Range("I110").Value = 521.8
$ws.Range("J110").Value = 599.75
$ws.Range("K110").Value = 521.8
$ws.Range("L110").Value = 599.75
$ws.Range("M110").Value = 1523.2
$ws.Range("N110").Value = -4689.75

$ws.Range("H116").Value = 1737.3715
$ws.Range("J116").Value = 2112.125
$ws.Range("L116").Value = 2112.125
$ws.Range("N116").Value = -6700.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1737.3715
$ws.Range("J3").Value = 2112.125
$ws.Range("L3").Value = 2112.125
$ws.Range("N3").Value = -2340.125

$ws.Range("H11").Value = 2949.75
$ws.Range("I11").Value = 3599.6667
$ws.Range("K11").Value = 3599.6667
$ws.Range("M11").Value = -3459.6667

$ws.Range("H20").Value = 1767.6364
$ws.Range("I20").Value = 2070.625
$ws.Range("J20").Value = 959.6667
$ws.Range("K20").Value = 2070.625
$ws.Range("L20").Value = 959.6667
$ws.Range("M20").Value = -1823.625
$ws.Range("N20").Value = -1453.6667

$ws.Range("H39").Value = 14000.0
$ws.Range("I39").Value = 14000.0
$ws.Range("K39").Value = 14000.0
$ws.Range("M39").Value = -13611.0

$ws.Range("H105").Value = 4356.7
$ws.Range("I105").Value = 5513.6
$ws.Range("K105").Value = 5513.6
$ws.Range("M105").Value = -3766.6

$ws.Range("H118").Value = 43650.0
$ws.Range("J118").Value = 43650.0
$ws.Range("L118").Value = 43650.0
$ws.Range("N118").Value = -46964.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11834.22
$ws.Range("I31").Value = 17049.52
$ws.Range("J31").Value = 3685.3125
$ws.Range("K31").Value = 17049.52
$ws.Range("L31").Value = 3685.3125
$ws.Range("M31").Value = -16754.52
$ws.Range("N31").Value = -4275.3125

$ws.Range("H34").Value = 11834.22
$ws.Range("I34").Value = 17049.52
$ws.Range("J34").Value = 3685.3125
$ws.Range("K34").Value = 17049.52
$ws.Range("L34").Value = 3685.3125
$ws.Range("M34").Value = -16847.52
$ws.Range("N34").Value = -4089.3125

$ws.Range("H122").Value = 1099.0322
$ws.Range("I122").Value = 938.05554
$ws.Range("J122").Value = 1321.9231
$ws.Range("K122").Value = 2814.16662
$ws.Range("L122").Value = 3965.7693
$ws.Range("M122").Value = -364.16662
$ws.Range("N122").Value = -8865.7693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3619.3572
$ws.Range("I2").Value = 6688.8
$ws.Range("J2").Value = 77.69231
$ws.Range("K2").Value = 40132.8
$ws.Range("L2").Value = 466.15386
$ws.Range("M2").Value = -40019.8
$ws.Range("N2").Value = -692.15386

$ws.Range("H5").Value = 1496.1538
$ws.Range("I5").Value = 1095.6364
$ws.Range("K5").Value = 3286.9092
$ws.Range("M5").Value = -3174.9092

$ws.Range("H23").Value = 1569.9
$ws.Range("J23").Value = 1577.6666
$ws.Range("L23").Value = 4732.9998
$ws.Range("N23").Value = -5202.9998

$ws.Range("H92").Value = 17857412.0
$ws.Range("I92").Value = 31250296.0
$ws.Range("J92").Value = 233.66667
$ws.Range("K92").Value = 93750888.0
$ws.Range("L92").Value = 701.00001
$ws.Range("M92").Value = -93749640.0
$ws.Range("N92").Value = -3197.00001

$ws.Range("H112").Value = 1304.0
$ws.Range("I112").Value = 630.0
$ws.Range("K112").Value = 1890.0
$ws.Range("M112").Value = -782.0

$ws.Range("H131").Value = 811.14
$ws.Range("J131").Value = 811.2525
$ws.Range("L131").Value = 2433.7575
$ws.Range("N131").Value = -12513.7575

$ws.Range("H135").Value = 1496.1538
$ws.Range("I135").Value = 1095.6364
$ws.Range("K135").Value = 9860.7276
$ws.Range("M135").Value = -7325.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3462.7
$ws.Range("I40").Value = 2637.0
$ws.Range("K40").Value = 2637.0
$ws.Range("M40").Value = -2501.0

$ws.Range("H82").Value = 1934.6786
$ws.Range("I82").Value = 1794.625
$ws.Range("J82").Value = 2775.0
$ws.Range("K82").Value = 1794.625
$ws.Range("L82").Value = 2775.0
$ws.Range("M82").Value = -1433.625
$ws.Range("N82").Value = -3497.0

$ws.Range("H85").Value = 1934.6786
$ws.Range("I85").Value = 1794.625
$ws.Range("J85").Value = 2775.0
$ws.Range("K85").Value = 1794.625
$ws.Range("L85").Value = 2775.0
$ws.Range("M85").Value = -546.625
$ws.Range("N85").Value = -5271.0

$ws.Range("H93").Value = 1261.8182
$ws.Range("I93").Value = 1247.6
$ws.Range("J93").Value = 1404.0
$ws.Range("K93").Value = 1247.6
$ws.Range("L93").Value = 1404.0
$ws.Range("M93").Value = 0.4000000000000909
$ws.Range("N93").Value = -3900.0

$ws.Range("H122").Value = 983186.1
$ws.Range("I122").Value = 2181380.2
$ws.Range("K122").Value = 6544140.600000001
$ws.Range("M122").Value = -6541690.600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 470.33334
$ws.Range("I100").Value = 469.0
$ws.Range("K100").Value = 938.0
$ws.Range("M100").Value = -397.0

$ws.Range("H122").Value = 1344.8334
$ws.Range("I122").Value = 1140.0
$ws.Range("J122").Value = 1549.6666
$ws.Range("K122").Value = 3420.0
$ws.Range("L122").Value = 4648.9998
$ws.Range("M122").Value = -970.0
$ws.Range("N122").Value = -9548.9998

$ws.Range("H132").Value = 829.675
$ws.Range("I132").Value = 568.5862
$ws.Range("J132").Value = 1518.0
$ws.Range("K132").Value = 1705.7586
$ws.Range("L132").Value = 4554.0
$ws.Range("M132").Value = 824.2414000000001
$ws.Range("N132").Value = -9614.0
